$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$wdStyleTypeParagraph = 1

# --- Introduce the new "TheoremStyleUpright" paragraph style -----------------
# It is a plain, un-italicised variant of "TheoremStyle" that the author now
# uses to control theorem styling explicitly (author control of theorem style
# classification / removal of italics).
$theoremUpright = $null
try {
    $theoremUpright = $d.Styles("TheoremStyleUpright")
} catch {
    $theoremUpright = $null
}
if (-not $theoremUpright) {
    $theoremUpright = $d.Styles.Add("TheoremStyleUpright", $wdStyleTypeParagraph)
}
$theoremUpright.BaseStyle = $d.Styles("TheoremStyle")
$theoremUpright.QuickStyle = $true

# --- Re-point the existing "Theorem" paragraphs at the new style -------------
foreach ($para in $d.Paragraphs) {
    if ($para.Style.NameLocal -eq "TheoremStyle") {
        $para.Style = $theoremUpright
    }
}
